# Strip the redundant "event_types/" directory prefix from the image
# column now that the loader resolves image paths relative to that
# folder itself (used by the new rspec image-testing coverage).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "event_types/"

for ($r = 2; $r -le 34; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value()
    if ($val -is [string] -and $val.StartsWith($prefix)) {
        $cell.Value = $val.Substring($prefix.Length)
    }
}

# Match Excel's auto row-height bump on the first edited data row.
$ws.Rows.Item(2).RowHeight = 16

$ws.Range("D2").Select()
